# Fix typo in column (task) names: the header row used strings like
# "TasP1".."TasP27" which should read "Task1".."Task27".
# A single Find & Replace over the whole sheet reproduces exactly what
# Excel does (and what the shared-strings table ends up looking like):
# every "TasP" occurrence becomes "Task".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("TasP", "Task", 2, 1, $false)
